# WIP: Last working copy (with surrounding circles)
# Adds a new Agent "A10" (with Knowledge personalization/link + a new linked
# Tool "T10") to the hierarchy tracked across the three sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("L12A")        # sheet1
$ws2 = $wb.Worksheets.Item("A2T")         # sheet2
$ws3 = $wb.Worksheets.Item("Components")  # sheet3

# --- L12A: fix "K" -> "Knowledge" label, then append new Agent A10 row ---
$ws1.Range("B5").Value = "Knowledge"
$ws1.Range("B6").Value = "Knowledge"
$ws1.Range("B7").Value = "Knowledge"
$ws1.Range("C7").Value = "A10"

$ws1.Columns.Item(2).ColumnWidth = 11.09

# --- A2T: append new Agent -> Tool mapping row ---
$ws2.Range("B6").Value = "A10"
$ws2.Range("C6").Value = "T10"

# --- Components: insert a new Agent row (A10) right after the existing
#     agents, pushing the Tool rows down by one ---
$ws3.Rows.Item(11).Insert()
$ws3.Range("B11").Value = "Agent"
$ws3.Range("C11").Value = "A10"
$ws3.Range("D11").Value = 1

# --- Selection / view bookkeeping (match the last-saved cursor positions) ---
# Order matters: the last sheet selected on is the one left as the active tab.
[void]$ws1.Range("C7").Select()
[void]$ws2.Range("C7").Select()
[void]$ws3.Range("E15").Select()
